# "Returned data table results to false"
# Sheet "Test Results" (4th sheet) has its CRUD-test-passed boolean columns
# (B:E, rows 2-24) reset back to FALSE, and the active selection moves to I18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")

$ws.Range("B2:E24").Value = $false

$ws.Activate()
$ws.Range("I18").Select()
